$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("<id>p150v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p150v_1</id>", 2)
